# CIV-6538: fix "must uploaded" -> "must upload" in the ADR / Digital Portal
# witness statement paragraph.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "must uploaded to the Digital Portal",  # FindText
    $true,                                   # MatchCase
    $false,                                  # MatchWholeWord
    $false,                                  # MatchWildcards
    $false,                                  # MatchSoundsLike
    $false,                                  # MatchAllWordForms
    $true,                                   # Forward
    1,                                        # Wrap (wdFindContinue)
    $false,                                  # Format
    "must upload to the Digital Portal",     # ReplaceWith
    2                                         # Replace (wdReplaceAll)
)
